$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.768.33'
$ws.Range('E2').Value = '  -1.68%  '
$ws.Range('D3').Value = '1.862.30'
$ws.Range('E3').Value = '  -2.75%  '
$ws.Range('E4').Value = '  -0.87%  '
$ws.Range('D5').Value = '''244.17'
$ws.Range('E5').Value = '  -3.85%  '
$ws.Range('D6').Value = '''0.674'
$ws.Range('E6').Value = '  -6.97%  '
$ws.Range('E7').Value = '  -0.90%  '
$ws.Range('D8').Value = '''41.60'
$ws.Range('E8').Value = '  +1.94%  '
$ws.Range('D9').Value = '''0.340'
$ws.Range('E9').Value = '  -4.50%  '
$ws.Range('E10').Value = '  -3.37%  '
$ws.Range('D11').Value = '''0.0966'
$ws.Range('E11').Value = '  -2.66%  '
$ws.Range('D12').Value = '''12.83'
$ws.Range('E12').Value = '  +1.11%  '
$ws.Range('D13').Value = '2.133.52'
$ws.Range('E13').Value = '  -2.77%  '
$ws.Range('D14').Value = '''0.707'
$ws.Range('E14').Value = '  -1.64%  '
$ws.Range('D15').Value = '1.876.63'
$ws.Range('E15').Value = '  -2.15%  '
$ws.Range('D16').Value = '''4.80'
$ws.Range('E16').Value = '  -2.80%  '
$ws.Range('D17').Value = '34.725.02'
$ws.Range('E17').Value = '  -1.87%  '
$ws.Range('D18').Value = '''72.01'
$ws.Range('E18').Value = '  -3.33%  '
$ws.Range('D19').Value = '0.0₃0808'
$ws.Range('E19').Value = '  -4.33%  '
$ws.Range('D20').Value = '''242.27'
$ws.Range('E20').Value = '  -0.78%  '
$ws.Range('D21').Value = '''12.51'
$ws.Range('E21').Value = '  -4.53%  '
$ws.Range('D22').Value = '''4.85'
$ws.Range('E22').Value = '  -4.68%  '
$ws.Range('E23').Value = '  -0.91%  '
$ws.Range('E24').Value = '  +4.82%  '
$ws.Range('D25').Value = '''2.12'
$ws.Range('E25').Value = '  -14.21%  '
$ws.Range('D26').Value = '''163.00'
$ws.Range('E26').Value = '  -2.27%  '
$ws.Range('D27').Value = '''8.30'
$ws.Range('E27').Value = '  -3.88%  '
$ws.Range('D28').Value = '''18.01'
$ws.Range('E29').Value = '  -5.88%  '
$ws.Range('D30').Value = '4.128.52'
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('E31').Value = '  +4.01%  '
$ws.Range('D32').Value = '''4.15'
$ws.Range('E32').Value = '  -5.38%  '
$ws.Range('D33').Value = '''0.0569'
$ws.Range('E33').Value = '  -1.99%  '
$ws.Range('E34').Value = '  -0.87%  '
$ws.Range('D35').Value = '''4.10'
$ws.Range('E35').Value = '  -3.15%  '
$ws.Range('D36').Value = '''0.823'
$ws.Range('E36').Value = '  -11.09%  '
$ws.Range('D37').Value = '''1.57'
$ws.Range('E37').Value = '  -21.33%  '
$ws.Range('D38').Value = '''1.94'
$ws.Range('E38').Value = '  -4.17%  '
$ws.Range('E39').Value = '  -0.23%  '
$ws.Range('D40').Value = '''16.87'
$ws.Range('E40').Value = '  -3.90%  '
$ws.Range('D41').Value = '''0.0662'
$ws.Range('E41').Value = '  +0.96%  '
$ws.Range('E42').Value = '  -4.03%  '
$ws.Range('D43').Value = '''1.06'
$ws.Range('E43').Value = '  -5.12%  '
$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D44').Value = '''0.0823'
$ws.Range('E44').Value = '  +11.38%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '1.278.90'
$ws.Range('E45').Value = '  -4.61%  '
$ws.Range('E47').Value = '  -1.13%  '
$ws.Range('E48').Value = '  -1.80%  '
$ws.Range('D49').Value = '''11.67'
$ws.Range('E49').Value = '  -1.44%  '
$ws.Range('D50').Value = '''6.23'
$ws.Range('E50').Value = '  -7.71%  '
$ws.Range('D51').Value = '''42.17'
$ws.Range('E51').Value = '  -6.65%  '
